$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 689 (shifts existing rows 689-739 down to 691-741)
$ws.Rows.Item(689).Insert()
$ws.Rows.Item(689).Insert()

# New row 689: Poroto verde, Magnum, Primera - Región Metropolitana
$ws.Cells.Item(689,1).Value = 9
$ws.Cells.Item(689,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(689,3).Value = "Metropolitana"
$ws.Cells.Item(689,4).Value = 45021
$ws.Cells.Item(689,5).Value = 13
$ws.Cells.Item(689,6).Value = 100112031
$ws.Cells.Item(689,7).Value = "Poroto verde"
$ws.Cells.Item(689,8).Value = "Magnum"
$ws.Cells.Item(689,9).Value = "Primera"
$ws.Cells.Item(689,10).Value = 40
$ws.Cells.Item(689,11).Value = 22000
$ws.Cells.Item(689,12).Value = 22000
$ws.Cells.Item(689,13).Value = 22000
$ws.Cells.Item(689,14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(689,15).Value = "Región Metropolitana"
$ws.Cells.Item(689,16).Value = 880
$ws.Cells.Item(689,17).Value = 25
$ws.Cells.Item(689,18).Value = "Hortaliza"

# New row 690: Poroto verde, Magnum, Segunda - Región Metropolitana
$ws.Cells.Item(690,1).Value = 9
$ws.Cells.Item(690,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(690,3).Value = "Metropolitana"
$ws.Cells.Item(690,4).Value = 45021
$ws.Cells.Item(690,5).Value = 13
$ws.Cells.Item(690,6).Value = 100112031
$ws.Cells.Item(690,7).Value = "Poroto verde"
$ws.Cells.Item(690,8).Value = "Magnum"
$ws.Cells.Item(690,9).Value = "Segunda"
$ws.Cells.Item(690,10).Value = 60
$ws.Cells.Item(690,11).Value = 23000
$ws.Cells.Item(690,12).Value = 23000
$ws.Cells.Item(690,13).Value = 23000
$ws.Cells.Item(690,14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(690,15).Value = "Región Metropolitana"
$ws.Cells.Item(690,16).Value = 920
$ws.Cells.Item(690,17).Value = 25
$ws.Cells.Item(690,18).Value = "Hortaliza"
